$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the cryptos price (D) and 1h-volume-change (E) columns.
# Force text format first so numeric-looking strings (e.g. "1.000",
# "30.403.84") are kept verbatim instead of being parsed as numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.403.84'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.08%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.873.19'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.64%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.10%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '235.78'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -1.07%  '

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.07%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4677'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.32%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2846'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.63%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06567'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.26%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.13'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +6.56%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07955'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +2.61%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '97.55'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -1.22%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.862.83'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -1.18%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.146'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.12%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6750'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.92%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '282.54'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -1.39%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '30.397.85'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.02%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '5.538'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +2.98%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.000'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.12%  '

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.12%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.116.07'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.66%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.000007304'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.17%  '

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.16%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.207'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.18%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.290'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.00%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '164.18'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -2.11%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.17'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +0.23%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.946'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -2.48%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.362'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -1.28%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.09716'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -1.28%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.445'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.92%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.478'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -1.12%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.116'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -1.95%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04693'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -0.51%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.116'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +1.51%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7057'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -0.84%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.714'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +0.01%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01862'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -0.86%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.334'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -5.67%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.544'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +0.49%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '73.48'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +0.85%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.946'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -1.35%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.8491'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -3.07%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.4197'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -0.27%  '

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.21%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '103.65'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -0.86%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '7.234'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -0.45%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.260'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -1.95%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '936.73'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -5.84%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '34.18'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -0.19%  '

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -2.79%  '
